$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Señales de control")

# Reorder the "Operaciones" table (G13:G19) to reflect the new
# instruction-fetch / instruction-decode stop module ordering.
$ws.Range("G13").Value = "SUB"
$ws.Range("G14").Value = "MOV"
$ws.Range("G15").Value = "MOVT"
$ws.Range("G16").Value = "AND"
$ws.Range("G18").Value = "EOR"
$ws.Range("G19").Value = "CMP"

# Update the active view: scrolled so A4 is the top-left visible cell,
# with the selection moved to G21.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G21").Select()
